{"js": "// Update the date line at the top of the document.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].getRange().insertText(\"2024-03-20 Wednesday\", Word.InsertLocation.replace);\n\n// Update the division problems in the table, cell by cell (positional,\n// since several before/after values repeat or collide with each other -\n// e.g. \"76\u00f79=\" is both a source and a target value - a single document-wide\n// search/replace would be unsafe).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// 0-based row indices of the five rows that hold the 5 problems each.\nconst newValues = {\n  0:  [\"93\u00f76=\", \"38\u00f73=\", \"36\u00f76=\", \"30\u00f73=\", \"52\u00f73=\"],\n  4:  [\"19\u00f74=\", \"58\u00f73=\", \"14\u00f73=\", \"16\u00f78=\", \"59\u00f74=\"],\n  8:  [\"43\u00f75=\", \"33\u00f75=\", \"31\u00f77=\", \"69\u00f76=\", \"76\u00f79=\"],\n  12: [\"42\u00f78=\", \"98\u00f77=\", \"71\u00f78=\", \"31\u00f73=\", \"93\u00f74=\"],\n  16: [\"32\u00f74=\", \"63\u00f76=\", \"94\u00f74=\", \"43\u00f78=\", \"39\u00f74=\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const row = parseInt(rowIndex, 10);\n  const vals = newValues[rowIndex];\n  for (let col = 0; col < vals.length; col++) {\n    const cell = table.getCell(row, col);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    cellParagraphs.items[0].getRange().insertText(vals[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line at the top of the document.\n$d.Content.Find.Execute(\"2024-03-19 Tuesday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-03-20 Wednesday\", 2)\n\n# Update the division problems in the table, cell by cell (positional,\n# since several before/after values repeat or collide with each other,\n# a single document-wide Find/Replace would be unsafe).\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"93\u00f76=\", \"38\u00f73=\", \"36\u00f76=\", \"30\u00f73=\", \"52\u00f73=\")\n    5  = @(\"19\u00f74=\", \"58\u00f73=\", \"14\u00f73=\", \"16\u00f78=\", \"59\u00f74=\")\n    9  = @(\"43\u00f75=\", \"33\u00f75=\", \"31\u00f77=\", \"69\u00f76=\", \"76\u00f79=\")\n    13 = @(\"42\u00f78=\", \"98\u00f77=\", \"71\u00f78=\", \"31\u00f73=\", \"93\u00f74=\")\n    17 = @(\"32\u00f74=\", \"63\u00f76=\", \"94\u00f74=\", \"43\u00f78=\", \"39\u00f74=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $vals = $newValues[$row]\n    for ($col = 1; $col -le $vals.Length; $col++) {\n        $cell = $t.Cell($row, $col)\n        $cellRange = $cell.Range\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        $cellRange.Text = $vals[$col - 1]\n    }\n}\n"}
